{"js": "// Hot fix: insert the \"{index}\" merge-field right after each\n// \"{#tipoDeElemento==...}\" block-opening tag, before the leading \"- Un/Una\"\n// text, for every element-type template block in the document body.\n//\n// Example:\n//   {#tipoDeElemento==\u201dno peritable\u201d}- Un (01) {descripcionElemento}, ...\n// becomes\n//   {#tipoDeElemento==\u201dno peritable\u201d}{index}- Un (01) {descripcionElemento}, ...\n\nconst body = context.document.body;\n\n// The curly quote characters differ between blocks in the source template\n// (some use RIGHT DOUBLE QUOTATION MARK \u201d U+201D on both sides, others use\n// LEFT \u201c U+201C / RIGHT \u201d U+201D), so each marker is listed explicitly\n// exactly as it appears in the document.\nconst markers = [\n  \"{#tipoDeElemento==\\u201Dno peritable\\u201D}\",\n  \"{#tipoDeElemento==\\u201Dsim\\u201D}\",\n  \"{#tipoDeElemento==\\u201Ddisco\\u201D}\",\n  \"{#tipoDeElemento==\\u201Cdvr\\u201D}\",\n  \"{#tipoDeElemento==\\u201Cgabinete\\u201D}\",\n  \"{#tipoDeElemento==\\u201Ctablet\\u201D}\",\n  \"{#tipoDeElemento==\\u201Dunidad de almacenamiento\\u201D}\",\n  \"{#tipoDeElemento==\\u201Cnotebook\\u201D}\",\n  \"{#tipoDeElemento==\\u201Ccelular\\u201D}\"\n];\n\nfor (const marker of markers) {\n  const results = body.search(marker, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"{index}\", Word.InsertLocation.after);\n  }\n  await context.sync();\n}\n", "ps1": "# Hot fix: insert the \"{index}\" merge-field right after each\n# \"{#tipoDeElemento==...}\" block-opening tag, before the leading \"- Un/Una\"\n# text, for every element-type template block in the document body.\n#\n# Example:\n#   {#tipoDeElemento==\u201dno peritable\u201d}- Un (01) {descripcionElemento}, ...\n# becomes\n#   {#tipoDeElemento==\u201dno peritable\u201d}{index}- Un (01) {descripcionElemento}, ...\n\n$d = $word.ActiveDocument\n\n# The curly quote characters differ between blocks in the source template\n# (some use RIGHT DOUBLE QUOTATION MARK \u201d U+201D on both sides, others use\n# LEFT \u201c U+201C / RIGHT \u201d U+201D), so each marker is built explicitly with\n# its exact characters as it appears in the document.\n$rdquo = [char]0x201D\n$ldquo = [char]0x201C\n\n$markers = @(\n    \"{#tipoDeElemento==\" + $rdquo + \"no peritable\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $rdquo + \"sim\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $rdquo + \"disco\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $ldquo + \"dvr\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $ldquo + \"gabinete\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $ldquo + \"tablet\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $rdquo + \"unidad de almacenamiento\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $ldquo + \"notebook\" + $rdquo + \"}\",\n    \"{#tipoDeElemento==\" + $ldquo + \"celular\" + $rdquo + \"}\"\n)\n\nforeach ($marker in $markers) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $newText = $marker + \"{index}\"\n    $find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
